$wb = $excel.ActiveWorkbook

# --- "Runs" sheet: insert a new row at the top for the merged run ---
$ws = $wb.Worksheets.Item("Runs")
$ws.Rows.Item(1).Insert()

$runsRow1 = @{
    1 = 0
    2 = 'Run076111114_final'
    3 = '0.697 (0.394)'
    4 = '0.688 (0.182)'
    5 = '0.701 (0.459)'
    6 = '0.787 (0.332)'
    7 = '0.663 (0.197)'
    8 = '0.844 (0.364)'
    9 = '0.428 (0.455)'
    10 = '0.586 (0.239)'
    11 = '0.387 (0.488)'
    12 = '0.628 (0.384)'
    13 = '0.660 (0.212)'
    14 = '0.598 (0.492)'
    15 = '0.631 (0.427)'
    16 = '0.646 (0.280)'
    17 = '0.623 (0.486)'
    18 = '0.178 (0.332)'
    19 = '0.636 (0.185)'
    20 = '0.075 (0.264)'
    21 = '0.810 (0.335)'
    22 = '0.632 (0.261)'
    23 = '0.870 (0.337)'
    24 = '0.864 (0.306)'
    25 = '0.537 (0.277)'
    26 = '0.917 (0.276)'
    27 = '0.531 (0.486)'
    28 = '0.484 (0.228)'
    29 = '0.534 (0.499)'
    30 = '0.733 (0.375)'
    31 = '0.661 (0.221)'
    32 = '0.766 (0.423)'
    33 = '0.792 (0.351)'
    34 = '0.627 (0.247)'
    35 = '0.844 (0.363)'
    36 = '0.427 (0.467)'
    37 = '0.578 (0.230)'
    38 = '0.402 (0.490)'
    39 = '0.544 (0.388)'
    40 = '0.639 (0.211)'
    41 = '0.440 (0.497)'
}
foreach ($c in $runsRow1.Keys) {
    $ws.Cells.Item(1, $c).Value = $runsRow1[$c]
}

# --- "GAN" sheet: populate row 1 with the merged-network stats ---
$gan = $wb.Worksheets.Item("GAN")
$ganRow1 = @{
    1 = 'Run006'
    2 = '19.345775132275133 (0.3463753781829967)'
    3 = '23.063606382978723 (1.3582629963740993)'
    4 = '20.60390476190476 (0.8939518105016319)'
    5 = '21.233763492063492 (0.5747396061533153)'
    6 = '0.12837566137566137 (0.017411381174361588)'
    7 = '0.40445478723404255 (0.05662784312527957)'
    8 = '0.15325396825396825 (0.018307221707849142)'
    9 = '0.21366031746031744 (0.03836058746689444)'
    10 = 0
    11 = '20.55998327759197 (1.523727171523039)'
    12 = '0.19685576923076925 (0.10171105620616296)'
}
foreach ($c in $ganRow1.Keys) {
    $gan.Cells.Item(1, $c).Value = $ganRow1[$c]
}
